# Apply "atualizacao dos dados da add" update:
# Row 57 (2025-05): retained_customers 143 -> 144, retention_rate recalculated
# Row 58 (2025-06): retained_customers 41 -> 49, prev_total_customers 220 -> 221, retention_rate recalculated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57
$ws.Range("B57").Value = 144
$ws.Range("D57").Value = 69.23076923076923

# Row 58
$ws.Range("B58").Value = 49
$ws.Range("C58").Value = 221
$ws.Range("D58").Value = 22.17194570135747
